$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. Drop the two trailing paragraphs ("…" and the blank one after it).
# The document must always keep a final paragraph mark, so merge away the
# "…" paragraph together with the empty paragraph that follows it, letting
# paragraph 5's own mark become the new end of story.
$p5End = $d.Paragraphs(5).Range.End
$p7End = $d.Paragraphs(7).Range.End
$d.Range($p5End, $p7End).Delete()

# --- 2. Paragraph 2 (still bordered): swap in the text that used to live in
# paragraph 3, keeping the original 3-run split.
$xmlP2 = @"
<w:p $wns>
  <w:pPr>
    <w:pBdr>
      <w:top w:val="single" w:sz="4" w:space="1" w:color="auto"/>
      <w:left w:val="single" w:sz="4" w:space="4" w:color="auto"/>
      <w:bottom w:val="single" w:sz="4" w:space="1" w:color="auto"/>
      <w:right w:val="single" w:sz="4" w:space="4" w:color="auto"/>
    </w:pBdr>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Add some comments about Version management </w:t></w:r>
  <w:r><w:t>outside this border</w:t></w:r>
  <w:r><w:t>, or just add some text so there is a change to this file.</w:t></w:r>
</w:p>
"@
$d.Paragraphs(2).Range.InsertXML($xmlP2)

# --- 3. Paragraph 3 (still bordered): swap in the text that used to live in
# paragraph 4, keeping the original 5-run split.
$xmlP3 = @"
<w:p $wns>
  <w:pPr>
    <w:pBdr>
      <w:top w:val="single" w:sz="4" w:space="1" w:color="auto"/>
      <w:left w:val="single" w:sz="4" w:space="4" w:color="auto"/>
      <w:bottom w:val="single" w:sz="4" w:space="1" w:color="auto"/>
      <w:right w:val="single" w:sz="4" w:space="4" w:color="auto"/>
    </w:pBdr>
  </w:pPr>
  <w:r><w:t xml:space="preserve">Remember that your GitHub user </w:t></w:r>
  <w:r><w:t>ID</w:t></w:r>
  <w:r><w:t xml:space="preserve"> must be submitted in you</w:t></w:r>
  <w:r><w:t>r</w:t></w:r>
  <w:r><w:t xml:space="preserve"> assignment report!</w:t></w:r>
</w:p>
"@
$d.Paragraphs(3).Range.InsertXML($xmlP3)

# --- 4. Paragraph 4 loses its border and becomes the lone ellipsis paragraph.
$xmlP4 = @"
<w:p $wns>
  <w:r><w:t>…</w:t></w:r>
</w:p>
"@
$d.Paragraphs(4).Range.InsertXML($xmlP4)

# --- 5. Paragraph 5 loses its border and becomes the new GitHub blurb, with
# a grammar-check proofing mark bracketing "all across".
$xmlP5 = @"
<w:p $wns>
  <w:r><w:t xml:space="preserve">Git Hub is an effective version management platform that allows team members </w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>all across</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> the world to work on a project.</w:t></w:r>
</w:p>
"@
$d.Paragraphs(5).Range.InsertXML($xmlP5)

"done"
